$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix up C37 timestamp precision
$ws.Cells.Item(37,3).Value = 45504.05657503472

# Append new price-log rows (38-59)
$ws.Cells.Item(38,1).Value = 'Cafea Organica House Roast, Exhale, boabe'
$ws.Cells.Item(38,2).Value = '155,00 Lei'
$ws.Cells.Item(38,3).Value = 45504.06547949074
$ws.Cells.Item(38,3).NumberFormat = 'yyyy-mm-dd h:mm:ss'

$ws.Cells.Item(39,1).Value = 'Klorane Urzica Sampon uscat 150 ml'
$ws.Cells.Item(39,2).Value = '35,05 Lei'
$ws.Cells.Item(39,3).Value = 45504.06548925926
$ws.Cells.Item(39,3).NumberFormat = 'yyyy-mm-dd h:mm:ss'

$ws.Cells.Item(40,1).Value = 'Crema pentru ochi Elmiplant Hyaluronic Gold, Femei, 15 ml'
$ws.Cells.Item(40,2).Value = '31,40 Lei'
$ws.Cells.Item(40,3).Value = 45504.06549770833
$ws.Cells.Item(40,3).NumberFormat = 'yyyy-mm-dd h:mm:ss'

$ws.Cells.Item(41,1).Value = 'Bautura de Ovaz Standard Minor Figures bax 6L'
$ws.Cells.Item(41,2).Value = '108,00 Lei'
$ws.Cells.Item(41,3).Value = 45504.06550571759
$ws.Cells.Item(41,3).NumberFormat = 'yyyy-mm-dd h:mm:ss'

$ws.Cells.Item(42,1).Value = 'Cafea Organica House Roast, Exhale, boabe'
$ws.Cells.Item(42,2).Value = '155,00 Lei'
$ws.Cells.Item(42,3).Value = 45504.06703585648
$ws.Cells.Item(42,3).NumberFormat = 'yyyy-mm-dd h:mm:ss'

$ws.Cells.Item(43,1).Value = 'Klorane Urzica Sampon uscat 150 ml'
$ws.Cells.Item(43,2).Value = '35,05 Lei'
$ws.Cells.Item(43,3).Value = 45504.06704974537
$ws.Cells.Item(43,3).NumberFormat = 'yyyy-mm-dd h:mm:ss'

$ws.Cells.Item(44,1).Value = 'Crema pentru ochi Elmiplant Hyaluronic Gold, Femei, 15 ml'
$ws.Cells.Item(44,2).Value = '31,40 Lei'
$ws.Cells.Item(44,3).Value = 45504.0670582176
$ws.Cells.Item(44,3).NumberFormat = 'yyyy-mm-dd h:mm:ss'

$ws.Cells.Item(45,1).Value = 'Bautura de Ovaz Standard Minor Figures bax 6L'
$ws.Cells.Item(45,2).Value = '108,00 Lei'
$ws.Cells.Item(45,3).Value = 45504.06706739583
$ws.Cells.Item(45,3).NumberFormat = 'yyyy-mm-dd h:mm:ss'

$ws.Cells.Item(46,1).Value = 'Cafea Organica House Roast, Exhale, boabe'
$ws.Cells.Item(46,2).Value = '155,00 Lei'
$ws.Cells.Item(46,3).Value = 45504.06756972223
$ws.Cells.Item(46,3).NumberFormat = 'yyyy-mm-dd h:mm:ss'

$ws.Cells.Item(47,1).Value = 'Klorane Urzica Sampon uscat 150 ml'
$ws.Cells.Item(47,2).Value = '35,05 Lei'
$ws.Cells.Item(47,3).Value = 45504.06757984954
$ws.Cells.Item(47,3).NumberFormat = 'yyyy-mm-dd h:mm:ss'

$ws.Cells.Item(48,1).Value = 'Crema pentru ochi Elmiplant Hyaluronic Gold, Femei, 15 ml'
$ws.Cells.Item(48,2).Value = '31,40 Lei'
$ws.Cells.Item(48,3).Value = 45504.06759200231
$ws.Cells.Item(48,3).NumberFormat = 'yyyy-mm-dd h:mm:ss'

$ws.Cells.Item(49,1).Value = 'Bautura de Ovaz Standard Minor Figures bax 6L'
$ws.Cells.Item(49,2).Value = '108,00 Lei'
$ws.Cells.Item(49,3).Value = 45504.06759958333
$ws.Cells.Item(49,3).NumberFormat = 'yyyy-mm-dd h:mm:ss'

$ws.Cells.Item(50,1).Value = 'Cafea Organica House Roast, Exhale, boabe'
$ws.Cells.Item(50,2).Value = '155,00 Lei'
$ws.Cells.Item(50,3).Value = 45504.06967789352
$ws.Cells.Item(50,3).NumberFormat = 'yyyy-mm-dd h:mm:ss'

$ws.Cells.Item(51,1).Value = 'Klorane Urzica Sampon uscat 150 ml'
$ws.Cells.Item(51,2).Value = '35,05 Lei'
$ws.Cells.Item(51,3).Value = 45504.06968862269
$ws.Cells.Item(51,3).NumberFormat = 'yyyy-mm-dd h:mm:ss'

$ws.Cells.Item(52,1).Value = 'Crema pentru ochi Elmiplant Hyaluronic Gold, Femei, 15 ml'
$ws.Cells.Item(52,2).Value = '31,40 Lei'
$ws.Cells.Item(52,3).Value = 45504.06969699074
$ws.Cells.Item(52,3).NumberFormat = 'yyyy-mm-dd h:mm:ss'

$ws.Cells.Item(53,1).Value = 'Bautura de Ovaz Standard Minor Figures bax 6L'
$ws.Cells.Item(53,2).Value = '108,00 Lei'
$ws.Cells.Item(53,3).Value = 45504.06970663195
$ws.Cells.Item(53,3).NumberFormat = 'yyyy-mm-dd h:mm:ss'

$ws.Cells.Item(54,1).Value = 'Prelata acoperire piscina, PVC, neagra, 366 cm, Bestway'
$ws.Cells.Item(54,2).Value = '149,80 Lei'
$ws.Cells.Item(54,3).Value = 45504.06972070602
$ws.Cells.Item(54,3).NumberFormat = 'yyyy-mm-dd h:mm:ss'

$ws.Cells.Item(55,1).Value = 'Cafea Organica House Roast, Exhale, boabe'
$ws.Cells.Item(55,2).Value = '155,00 Lei'
$ws.Cells.Item(55,3).Value = 45504.07166885417
$ws.Cells.Item(55,3).NumberFormat = 'yyyy-mm-dd h:mm:ss'

$ws.Cells.Item(56,1).Value = 'Klorane Urzica Sampon uscat 150 ml'
$ws.Cells.Item(56,2).Value = '35,05 Lei'
$ws.Cells.Item(56,3).Value = 45504.07167846065
$ws.Cells.Item(56,3).NumberFormat = 'yyyy-mm-dd h:mm:ss'

$ws.Cells.Item(57,1).Value = 'Crema pentru ochi Elmiplant Hyaluronic Gold, Femei, 15 ml'
$ws.Cells.Item(57,2).Value = '31,40 Lei'
$ws.Cells.Item(57,3).Value = 45504.07168659722
$ws.Cells.Item(57,3).NumberFormat = 'yyyy-mm-dd h:mm:ss'

$ws.Cells.Item(58,1).Value = 'Bautura de Ovaz Standard Minor Figures bax 6L'
$ws.Cells.Item(58,2).Value = '108,00 Lei'
$ws.Cells.Item(58,3).Value = 45504.07169497685
$ws.Cells.Item(58,3).NumberFormat = 'yyyy-mm-dd h:mm:ss'

$ws.Cells.Item(59,1).Value = 'Prelata acoperire piscina, PVC, neagra, 366 cm, Bestway'
$ws.Cells.Item(59,2).Value = '149,80 Lei'
$ws.Cells.Item(59,3).Value = 45504.07170875659
$ws.Cells.Item(59,3).NumberFormat = 'yyyy-mm-dd h:mm:ss'